$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking strings are preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.296.12"
$ws.Range("E2").Value = "  -3.26%  "

$ws.Range("D3").Value = "1.734.89"
$ws.Range("E3").Value = "  -3.85%  "

$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "323.18"
$ws.Range("E5").Value = "  -4.70%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.27%  "

$ws.Range("D7").Value = "0.4254"
$ws.Range("E7").Value = "  -13.46%  "

$ws.Range("D8").Value = "0.3609"
$ws.Range("E8").Value = "  -3.27%  "

$ws.Range("D9").Value = "44.91"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("E10").Value = "  -3.00%  "

$ws.Range("D11").Value = "0.07368"
$ws.Range("E11").Value = "  -5.15%  "

$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").Value = "21.50"
$ws.Range("E13").Value = "  -5.14%  "

$ws.Range("D14").Value = "6.064"
$ws.Range("E14").Value = "  -4.34%  "

$ws.Range("D15").Value = "7.172"
$ws.Range("E15").Value = "  -2.26%  "

$ws.Range("D16").Value = "1.736.06"
$ws.Range("E16").Value = "  -3.25%  "

$ws.Range("D17").Value = "0.00001058"
$ws.Range("E17").Value = "  -3.86%  "

$ws.Range("D18").Value = "84.53"
$ws.Range("E18").Value = "  +2.54%  "

$ws.Range("D19").Value = "0.05950"
$ws.Range("E19").Value = "  -11.74%  "

$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").Value = "16.78"
$ws.Range("E21").Value = "  -3.84%  "

$ws.Range("D22").Value = "6.011"
$ws.Range("E22").Value = "  -6.62%  "

$ws.Range("D23").Value = "27.337.31"
$ws.Range("E23").Value = "  -3.07%  "

$ws.Range("D24").Value = "11.26"
$ws.Range("E24").Value = "  -6.38%  "

$ws.Range("D25").Value = "2.398"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "19.87"
$ws.Range("E26").Value = "  -5.43%  "

$ws.Range("D27").Value = "2.331"
$ws.Range("E27").Value = "  -3.67%  "

$ws.Range("D28").Value = "148.58"
$ws.Range("E28").Value = "  -1.91%  "

$ws.Range("D29").Value = "1.937.06"
$ws.Range("E29").Value = "  -3.34%  "

$ws.Range("D30").Value = "1.249"
$ws.Range("E30").Value = "  -2.16%  "

$ws.Range("D31").Value = "125.78"
$ws.Range("E31").Value = "  -6.58%  "

$ws.Range("D32").Value = "3.715"
$ws.Range("E32").Value = "  -8.13%  "

$ws.Range("D33").Value = "0.09002"
$ws.Range("E33").Value = "  -9.29%  "

$ws.Range("D34").Value = "5.543"
$ws.Range("E34").Value = "  -7.11%  "

$ws.Range("D35").Value = "12.33"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("D36").Value = "0.2158"
$ws.Range("E36").Value = "  -3.06%  "

$ws.Range("D37").Value = "0.02268"
$ws.Range("E37").Value = "  -5.21%  "

$ws.Range("D38").Value = "0.06086"
$ws.Range("E38").Value = "  -4.78%  "

$ws.Range("D39").Value = "0.6404"
$ws.Range("E39").Value = "  -4.87%  "

$ws.Range("D40").Value = "4.985"
$ws.Range("E40").Value = "  -5.16%  "

$ws.Range("D41").Value = "1.179"
$ws.Range("E41").Value = "  -3.78%  "

$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("D43").Value = "1.412"
$ws.Range("E43").Value = "  -4.96%  "

$ws.Range("D44").Value = "7.799"
$ws.Range("E44").Value = "  -4.30%  "

$ws.Range("D45").Value = "13.38"
$ws.Range("E45").Value = "  -5.73%  "

$ws.Range("D46").Value = "3.743"
$ws.Range("E46").Value = "  -3.34%  "

$ws.Range("D47").Value = "0.5831"
$ws.Range("E47").Value = "  -5.70%  "

$ws.Range("D48").Value = "124.18"
$ws.Range("E48").Value = "  -4.01%  "

$ws.Range("D49").Value = "1.927"
$ws.Range("E49").Value = "  -6.58%  "

$ws.Range("D50").Value = "0.06802"
$ws.Range("E50").Value = "  -4.62%  "

$ws.Range("D51").Value = "1.091"
$ws.Range("E51").Value = "  -7.38%  "
